$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-12-05 Friday", $true, $false, $false, $false,
                         $false, $true, 1, $false, "2025-12-06 Saturday", 2) | Out-Null

# Update the multiplication problems in the table, cell by cell (some
# values repeat, e.g. "53x63=", so Find/Replace across the whole table
# would be ambiguous; direct cell addressing guarantees the correct
# occurrence is updated).
$t = $d.Tables.Item(1)

$updates = @{
    1  = @("74×87=", "20×64=", "24×81=", "93×11=", "51×83=")
    5  = @("92×66=", "60×13=", "67×37=", "97×87=", "99×79=")
    10 = @("90×34=", "44×45=", "11×92=", "27×78=", "79×83=")
    15 = @("95×79=", "18×22=", "36×75=", "83×70=", "41×26=")
    20 = @("21×77=", "58×79=", "35×89=", "49×65=", "18×56=")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    for ($col = 1; $col -le $values.Count; $col++) {
        $t.Cell($row, $col).Range.Text = $values[$col - 1]
    }
}
